$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 48461.9
$ws.Range("J134").Value = 48461.9
$ws.Range("L134").Value = 48461.9
$ws.Range("N134").Value = -58601.9
$ws.Range("H137").Value = 514685.38
$ws.Range("J137").Value = 2709.5178
$ws.Range("L137").Value = 8128.553400000001
$ws.Range("N137").Value = -13228.5534

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2729.9546
$ws.Range("I74").Value = 2248.2144
$ws.Range("J74").Value = 3573
$ws.Range("K74").Value = 2248.2144
$ws.Range("L74").Value = 3573
$ws.Range("M74").Value = -1374.2144
$ws.Range("N74").Value = -5321
$ws.Range("H77").Value = 2729.9546
$ws.Range("I77").Value = 2248.2144
$ws.Range("J77").Value = 3573
$ws.Range("K77").Value = 11241.072
$ws.Range("L77").Value = 17865
$ws.Range("M77").Value = -6873.072
$ws.Range("N77").Value = -26601
$ws.Range("H110").Value = 1011
$ws.Range("I110").Value = 894
$ws.Range("J110").Value = 1362
$ws.Range("K110").Value = 894
$ws.Range("L110").Value = 1362
$ws.Range("M110").Value = 1151
$ws.Range("N110").Value = -5452
$ws.Range("H122").Value = 3534.2222
$ws.Range("I122").Value = 3211.8928
$ws.Range("K122").Value = 9635.678400000001
$ws.Range("M122").Value = -7185.678400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1464.421
$ws.Range("I107").Value = 1492.0714
$ws.Range("J107").Value = 1387
$ws.Range("K107").Value = 1492.0714
$ws.Range("L107").Value = 1387
$ws.Range("M107").Value = 427.9286
$ws.Range("N107").Value = -5227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4446151.5
$ws.Range("I16").Value = 8548445
$ws.Range("J16").Value = 1999.9166
$ws.Range("K16").Value = 8548445
$ws.Range("L16").Value = 1999.9166
$ws.Range("M16").Value = -8548158
$ws.Range("N16").Value = -2573.9166
$ws.Range("H31").Value = 2144.197
$ws.Range("I31").Value = 942.5
$ws.Range("J31").Value = 2925.3
$ws.Range("K31").Value = 942.5
$ws.Range("L31").Value = 2925.3
$ws.Range("M31").Value = -647.5
$ws.Range("N31").Value = -3515.3
$ws.Range("H34").Value = 2144.197
$ws.Range("I34").Value = 942.5
$ws.Range("J34").Value = 2925.3
$ws.Range("K34").Value = 942.5
$ws.Range("L34").Value = 2925.3
$ws.Range("M34").Value = -740.5
$ws.Range("N34").Value = -3329.3
$ws.Range("H94").Value = 1113.96
$ws.Range("I94").Value = 517.7143
$ws.Range("J94").Value = 1345.8334
$ws.Range("K94").Value = 517.7143
$ws.Range("L94").Value = 1345.8334
$ws.Range("M94").Value = -66.71429999999998
$ws.Range("N94").Value = -2247.8334
$ws.Range("H107").Value = 697.875
$ws.Range("I107").Value = 395.27585
$ws.Range("J107").Value = 1495.6364
$ws.Range("K107").Value = 395.27585
$ws.Range("L107").Value = 1495.6364
$ws.Range("M107").Value = 1524.72415
$ws.Range("N107").Value = -5335.6364
$ws.Range("H113").Value = 4446151.5
$ws.Range("I113").Value = 8548445
$ws.Range("J113").Value = 1999.9166
$ws.Range("K113").Value = 8548445
$ws.Range("L113").Value = 1999.9166
$ws.Range("M113").Value = -8546275
$ws.Range("N113").Value = -6339.9166
$ws.Range("H125").Value = 35163
$ws.Range("J125").Value = 35163
$ws.Range("L125").Value = 35163
$ws.Range("N125").Value = -40083
$ws.Range("H132").Value = 2093.1333
$ws.Range("I132").Value = 1382.72
$ws.Range("J132").Value = 5645.2
$ws.Range("K132").Value = 4148.16
$ws.Range("L132").Value = 16935.6
$ws.Range("M132").Value = -1618.16
$ws.Range("N132").Value = -21995.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1539.491
$ws.Range("I68").Value = 1155.6
$ws.Range("J68").Value = 1683.45
$ws.Range("K68").Value = 3466.8
$ws.Range("L68").Value = 5050.35
$ws.Range("M68").Value = -2655.8
$ws.Range("N68").Value = -6672.35
$ws.Range("H71").Value = 1539.491
$ws.Range("I71").Value = 1155.6
$ws.Range("J71").Value = 1683.45
$ws.Range("K71").Value = 10400.4
$ws.Range("L71").Value = 15151.05
$ws.Range("M71").Value = -6344.4
$ws.Range("N71").Value = -23263.05
$ws.Range("H103").Value = 1012.5
$ws.Range("I103").Value = 1012.5
$ws.Range("K103").Value = 3037.5
$ws.Range("M103").Value = -2158.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 12346522
$ws.Range("I107").Value = 533.3333
$ws.Range("J107").Value = 18519516
$ws.Range("K107").Value = 533.3333
$ws.Range("L107").Value = 18519516
$ws.Range("M107").Value = 1386.6667
$ws.Range("N107").Value = -18523356
$ws.Range("H109").Value = 10028
$ws.Range("J109").Value = 10028
$ws.Range("L109").Value = 10028
$ws.Range("N109").Value = -12108
$ws.Range("H122").Value = 3944.1875
$ws.Range("I122").Value = 3413.1428
$ws.Range("J122").Value = 4357.222
$ws.Range("K122").Value = 10239.4284
$ws.Range("L122").Value = 13071.666
$ws.Range("M122").Value = -7789.428400000001
$ws.Range("N122").Value = -17971.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3042.7036
$ws.Range("I7").Value = 1309.5333
$ws.Range("J7").Value = 5209.1665
$ws.Range("K7").Value = 1309.5333
$ws.Range("L7").Value = 5209.1665
$ws.Range("M7").Value = -1197.5333
$ws.Range("N7").Value = -5433.1665
$ws.Range("H40").Value = 6285.846
$ws.Range("I40").Value = 5610.5454
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 5610.5454
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -5474.5454
$ws.Range("N40").Value = -10272
$ws.Range("H126").Value = 3042.7036
$ws.Range("I126").Value = 1309.5333
$ws.Range("J126").Value = 5209.1665
$ws.Range("K126").Value = 3928.5999
$ws.Range("L126").Value = 15627.4995
$ws.Range("M126").Value = -1458.5999
$ws.Range("N126").Value = -20567.4995
$ws.Range("H137").Value = 48360
$ws.Range("J137").Value = 48360
$ws.Range("L137").Value = 48360
$ws.Range("N137").Value = -58560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4799.409
$ws.Range("I122").Value = 2783.4285
$ws.Range("J122").Value = 5740.2
$ws.Range("K122").Value = 8350.2855
$ws.Range("L122").Value = 17220.6
$ws.Range("M122").Value = -5900.2855
$ws.Range("N122").Value = -22120.6
